# Lagde et filter for minimum inntekt for utbetalingstidslinjer
# Insert a new first data row (row 2) with a zero-baseline data point
# (date 2019-09-20 / 43728) so the "Unit Test Counts" chart shows the
# count history starting from zero, and shift the rest of the table
# down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row above the current row 2, pushing rows 2:7 down
#        to rows 3:8 (formulas/refs on the existing rows auto-adjust). ---
$ws.Rows("2:2").Insert()

# --- 2. Populate the new row 2 with the baseline data point. ---
$ws.Range("A2").Value = 43728
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0

# --- 3. Re-enter the D column "total" formula as shared-formula ranges
#        matching how Excel itself re-groups a shared formula that an
#        inserted row has split in two: the original block (rows that
#        already existed) keeps going as one group, and the newly
#        inserted row + its immediate neighbour form a second group. ---
$ws.Range("D2:D3").Formula = "=C2+B2"
$ws.Range("D4:D8").Formula = "=C4+B4"

# --- 4. The inserted row copied row 1's (bold header) style onto A2:D2;
#        clear that back to the default, then restore the date number
#        format on A2 so it matches the other date cells in column A. ---
$ws.Range("A2:D2").ClearFormats()
$ws.Range("A2").NumberFormat = "d-mmm"

# --- 5. Move the active selection (matches the saved workbook view). ---
$ws.Range("D10").Select() | Out-Null

# --- 6. Update both chart series so they include the new row 8 (the
#        table now runs A2:D8 instead of A2:D7). ---
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$series = $chart.SeriesCollection()

$s1 = $series.Item(1)
$s1.Formula = '=SERIES("Unit Test Counts",Sheet1!$A$2:$A$8,Sheet1!$B$2:$B$8,1)'

$s2 = $series.Item(2)
$s2.Formula = '=SERIES(Sheet1!$C$1,Sheet1!$A$2:$A$8,Sheet1!$C$2:$C$8,2)'

# --- 7. Grow the chart by one row's worth of height so its anchor keeps
#        covering the same number of data rows visually (the bottom
#        edge moves down from row 14 to row 15, offset unchanged). ---
$co.Height = $co.Height() + $ws.Rows.Item(2).RowHeight()
